$wb = $excel.ActiveWorkbook

# Rename the second sheet ("Include from ") to "Include from Simplified Modif"
$wsInclude = $wb.Worksheets.Item("Include from ")
$wsInclude.Name = "Include from Simplified Modif"

# Metadata sheet updates
$wsMeta = $wb.Worksheets.Item("Metadata")

# Version: 1.0.1 -> 0.0.0
$wsMeta.Range("B3").Value = "0.0.0"

# Title
$wsMeta.Range("B5").Value = "Simplified modified Ranking Scale questionnaire (SMRSq)"

# Experimental value: (blank) -> "false" stored as TEXT (not boolean).
# A bare Value = "false" auto-coerces to a Boolean cell, so enter it with a
# leading apostrophe (forces text) and then restore the normal cell format
# (copied from its neighbour A7) so the text-entry quote-prefix styling
# doesn't linger on the cell.
$cellB7 = $wsMeta.Range("B7")
$cellB7.Value = "'false"
$wsMeta.Range("A7").Copy()
$cellB7.PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Date
$wsMeta.Range("B8").Value = "2024-01-11T13:00:00-03:00"

# Description
$wsMeta.Range("B12").Value = "ValueSet that defines the response values for the simplified modified Ranking Scale questionnaire."

# Include sheet updates
# System URI value
$wsInclude.Range("B10").Value = "https://molic-avc.gabriellesantosleandro.com/CodeSystem/SMRSqCS"
